# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, bordered, centered) from the
# last header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Header labels for the new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every player row (2-51): 81 wins, 81 losses, 0 ties.
$ws.Range("AD2:AE51").Value = 81
$ws.Range("AF2:AF51").Value = 0
